$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Periodo Mora" column (E16:E22) values with the new period list,
# entered in reverse-chronological (newest-first) order. This removes the
# previous account-statement periods and adds the new ones, per the commit
# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos".
$periodos = @("2402", "2401", "2312", "2311", "2310", "2309", "2308")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
